# Apply the "Add files via upload" edit: re-label / reorder the header
# columns, fill in sample employee rows 2-4, add a right-aligned style used
# by the CONTRACT PERIOD column, resize a few columns, and move the active
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row (row 1) - re-label / move columns.
#    A1:D1 (ROLL NUMBER, EMPLOYEE NAME, GENDER, DOB) stay as-is.
# ---------------------------------------------------------------------
$ws.Range("E1").Value = "VISA EXPIRY"
$ws.Range("F1").Value = "CONTRACT PERIOD"
$ws.Range("G1").Value = "WORKING HOURS"
$ws.Range("H1").Value = "OVERTIME"
$ws.Range("I1").Value = "SALARY"
$ws.Range("J1").Value = "PERFORMANCE"
$ws.Range("K1").Value = "MANAGER FEEDBACK "
$ws.Range("L1").Value = "MANAGER RATING OUT OF 10"
$ws.Range("M1").Value = "OUTPUT "

# ---------------------------------------------------------------------
# 2. Fill sample data for the first three employees (rows 2-4).
# ---------------------------------------------------------------------
$ws.Range("B2").Value = "AHAMED"
$ws.Range("C2").Value = "MALE"
$ws.Range("D2").Value = 1988
$ws.Range("E2").Value = 2025
$ws.Range("F2").Value = "2 YRS"
$ws.Range("G2").Value = 8
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 3500
$ws.Range("J2").Value = "GOOD"
$ws.Range("K2").Value = "HARDWORKER"
$ws.Range("L2").Value = 8
$ws.Range("M2").Value = "NOT RESIGN "

$ws.Range("B3").Value = "ARUN KUMAR"
$ws.Range("C3").Value = "MALE"
$ws.Range("D3").Value = 1988
$ws.Range("E3").Value = 2025
$ws.Range("F3").Value = "2 YRS"
$ws.Range("G3").Value = 8
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 3500
$ws.Range("J3").Value = "NOT GOOD"
$ws.Range("K3").Value = "OFFICE TIMING NOT KEEPING AND TASK DOING SLOWLY"
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = "RESIGN"

$ws.Range("B4").Value = "ASHIK"
$ws.Range("C4").Value = "MALE"
$ws.Range("D4").Value = 1988
$ws.Range("E4").Value = 2025
$ws.Range("F4").Value = "2 YRS"
$ws.Range("G4").Value = 8
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 4000
$ws.Range("J4").Value = "GOOD"
$ws.Range("K4").Value = "HARDWORKER"
$ws.Range("L4").Value = 9
$ws.Range("M4").Value = "NOT RESIGN "

# ---------------------------------------------------------------------
# 3. New right-aligned style, used by the CONTRACT PERIOD column (F2:F16)
#    and (oddly, matching the source workbook) G3.
# ---------------------------------------------------------------------
$ws.Range("F2:F16").HorizontalAlignment = -4152   # xlRight
$ws.Range("G3").HorizontalAlignment = -4152       # xlRight

# ---------------------------------------------------------------------
# 4. Column width adjustments to roughly match the new content.
#    (Values are chosen so the stored width lands as close as possible
#    to the target widths.)  Column I no longer needs a custom "best
#    fit" width once it only holds numeric salaries, so it is set back
#    close to the workbook's standard width.
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 11.333333333333334    # E  -> ~12.14
$ws.Columns.Item(6).ColumnWidth = 17                     # F  -> ~17.86
$ws.Columns.Item(7).ColumnWidth = 15.666666666666666     # G  -> ~16.57
$ws.Columns.Item(8).ColumnWidth = 9.333333333333334      # H  -> ~10.14
$ws.Columns.Item(9).ColumnWidth = 7.666666666666667      # I  -> standard width
$ws.Columns.Item(10).ColumnWidth = 13.666666666666666    # J  -> ~14.57
$ws.Columns.Item(11).ColumnWidth = 50.166666666666664    # K  -> 51
$ws.Columns.Item(12).ColumnWidth = 26.5                  # L  -> ~27.29

# ---------------------------------------------------------------------
# 5. Move the active selection to N4 (matches the saved sheet view).
# ---------------------------------------------------------------------
$ws.Range("N4").Select()
